# Automatiza la carga de los archivos por jornada al archivo principal.
# Minutos, Goles, Asistencias
#
# 1) Columna "Capitan" (AN) para las filas de jugadores (2-21): se limpia
#    el valor "No"/"Si" que quedó de una carga anterior.
# 2) Columna "Suplente" (F) para las filas 13-21: corrige el texto "Si"
#    (sin tilde) a "Sí" (con tilde).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Limpia la columna Capitan (AN) para todas las filas de datos (2 a 21).
$ws.Range("AN2:AN21").ClearContents()

# Corrige la tilde de "Si" a "Sí" en la columna Suplente (F) de las filas 13 a 21.
for ($row = 13; $row -le 21; $row++) {
    $ws.Cells.Item($row, 6).Value = "Sí"
}
